# ============================================================================
# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" positioned right before "总计", built by
#    copying the "2021-Q4" sheet (same column layout/styles: bold bordered
#    header row + bold bordered index column A) and then overwriting its
#    20 data rows with the 2022-Q1 fund-holding figures.
# 2. Insert a new top data row into "总计" for the 2022-Q1 summary figures,
#    shifting the existing 2021-Q4 / 2021-Q3 / 2021-Q2 rows down and
#    renumbering the index column.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: build the "2022-Q1" sheet from the "2021-Q4" template
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$zjBeforeCopy = $wb.Worksheets.Item("总计")

# Copy "2021-Q4" to just before "总计"; the copy keeps the header/index
# column styling (border + bold, style index shared with other sheets) and
# all 20 data rows already sized correctly (A1:H21).
$template.Copy($zjBeforeCopy)

# The copy lands immediately before "总计" and Excel names it "2021-Q4 (2)".
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# Capture the existing (correct) per-column style of the data block before
# touching any values, so we can strip the "quote-prefix" style that typing
# a leading-apostrophe string produces (Excel COM otherwise coerces numeric
# looking strings like "004702" or "30.93" into real numbers, which is not
# what the source data uses - these columns are text).
$dataStyle = $new.Range("B2:G21").Style

$new.Range("A2").Value = 0
$new.Range("B2").Value = "'004702"
$new.Range("C2").Value = "'南方金融主题灵活配置混合"
$new.Range("D2").Value = "'30.93"
$new.Range("E2").Value = "'89.81"
$new.Range("F2").Value = "'4.06"
$new.Range("G2").Value = "'1.2558"
$new.Range("H2").Value = 10
$new.Range("A3").Value = 1
$new.Range("B3").Value = "'005662"
$new.Range("C3").Value = "'嘉实金融精选股票A"
$new.Range("D3").Value = "'10.58"
$new.Range("E3").Value = "'90.89"
$new.Range("F3").Value = "'8.08"
$new.Range("G3").Value = "'0.8549"
$new.Range("H3").Value = 4
$new.Range("A4").Value = 2
$new.Range("B4").Value = "'000061"
$new.Range("C4").Value = "'华夏盛世混合"
$new.Range("D4").Value = "'14.90"
$new.Range("E4").Value = "'86.80"
$new.Range("F4").Value = "'2.20"
$new.Range("G4").Value = "'0.3278"
$new.Range("H4").Value = 8
$new.Range("A5").Value = 3
$new.Range("B5").Value = "'370024"
$new.Range("C5").Value = "'上投摩根核心优选混合"
$new.Range("D5").Value = "'13.74"
$new.Range("E5").Value = "'85.27"
$new.Range("F5").Value = "'1.93"
$new.Range("G5").Value = "'0.2652"
$new.Range("H5").Value = 10
$new.Range("A6").Value = 4
$new.Range("B6").Value = "'377530"
$new.Range("C6").Value = "'上投摩根行业轮动混合A"
$new.Range("D6").Value = "'10.38"
$new.Range("E6").Value = "'86.32"
$new.Range("F6").Value = "'2.44"
$new.Range("G6").Value = "'0.2533"
$new.Range("H6").Value = 9
$new.Range("A7").Value = 5
$new.Range("B7").Value = "'004686"
$new.Range("C7").Value = "'华夏研究精选股票"
$new.Range("D7").Value = "'4.76"
$new.Range("E7").Value = "'92.92"
$new.Range("F7").Value = "'5.31"
$new.Range("G7").Value = "'0.2528"
$new.Range("H7").Value = 2
$new.Range("A8").Value = 6
$new.Range("B8").Value = "'000082"
$new.Range("C8").Value = "'嘉实研究阿尔法股票"
$new.Range("D8").Value = "'7.51"
$new.Range("E8").Value = "'89.53"
$new.Range("F8").Value = "'3.28"
$new.Range("G8").Value = "'0.2463"
$new.Range("H8").Value = 1
$new.Range("A9").Value = 7
$new.Range("B9").Value = "'005663"
$new.Range("C9").Value = "'嘉实金融精选股票C"
$new.Range("D9").Value = "'3.00"
$new.Range("E9").Value = "'90.89"
$new.Range("F9").Value = "'8.08"
$new.Range("G9").Value = "'0.2424"
$new.Range("H9").Value = 4
$new.Range("A10").Value = 8
$new.Range("B10").Value = "'011911"
$new.Range("C10").Value = "'华夏消费优选混合型证券投资基金A"
$new.Range("D10").Value = "'7.18"
$new.Range("E10").Value = "'82.18"
$new.Range("F10").Value = "'2.64"
$new.Range("G10").Value = "'0.1896"
$new.Range("H10").Value = 8
$new.Range("A11").Value = 9
$new.Range("B11").Value = "'960006"
$new.Range("C11").Value = "'上投摩根行业轮动混合H"
$new.Range("D11").Value = "'2.00"
$new.Range("E11").Value = "'86.32"
$new.Range("F11").Value = "'2.44"
$new.Range("G11").Value = "'0.0488"
$new.Range("H11").Value = 9
$new.Range("A12").Value = 10
$new.Range("B12").Value = "'012244"
$new.Range("C12").Value = "'广发金融地产精选股票型发起式证券投资基金A"
$new.Range("D12").Value = "'0.98"
$new.Range("E12").Value = "'86.81"
$new.Range("F12").Value = "'4.80"
$new.Range("G12").Value = "'0.0470"
$new.Range("H12").Value = 2
$new.Range("A13").Value = 11
$new.Range("B13").Value = "'001254"
$new.Range("C13").Value = "'泰达宏利新起点灵活配置混合A"
$new.Range("D13").Value = "'1.88"
$new.Range("E13").Value = "'30.41"
$new.Range("F13").Value = "'1.76"
$new.Range("G13").Value = "'0.0331"
$new.Range("H13").Value = 5
$new.Range("A14").Value = 12
$new.Range("B14").Value = "'012245"
$new.Range("C14").Value = "'广发金融地产精选股票型发起式证券投资基金C"
$new.Range("D14").Value = "'0.63"
$new.Range("E14").Value = "'86.81"
$new.Range("F14").Value = "'4.80"
$new.Range("G14").Value = "'0.0302"
$new.Range("H14").Value = 2
$new.Range("A15").Value = 13
$new.Range("B15").Value = "'160722"
$new.Range("C15").Value = "'嘉实惠泽灵活配置混合(LOF)"
$new.Range("D15").Value = "'0.67"
$new.Range("E15").Value = "'92.76"
$new.Range("F15").Value = "'2.93"
$new.Range("G15").Value = "'0.0196"
$new.Range("H15").Value = 8
$new.Range("A16").Value = 14
$new.Range("B16").Value = "'011912"
$new.Range("C16").Value = "'华夏消费优选混合型证券投资基金C"
$new.Range("D16").Value = "'0.44"
$new.Range("E16").Value = "'82.18"
$new.Range("F16").Value = "'2.64"
$new.Range("G16").Value = "'0.0116"
$new.Range("H16").Value = 8
$new.Range("A17").Value = 15
$new.Range("B17").Value = "'162211"
$new.Range("C17").Value = "'泰达宏利品质生活混合"
$new.Range("D17").Value = "'0.12"
$new.Range("E17").Value = "'76.13"
$new.Range("F17").Value = "'4.23"
$new.Range("G17").Value = "'0.0051"
$new.Range("H17").Value = 7
$new.Range("A18").Value = 16
$new.Range("B18").Value = "'002313"
$new.Range("C18").Value = "'泰达宏利新起点灵活配置混合B"
$new.Range("D18").Value = "'0.16"
$new.Range("E18").Value = "'30.41"
$new.Range("F18").Value = "'1.76"
$new.Range("G18").Value = "'0.0028"
$new.Range("H18").Value = 5
$new.Range("A19").Value = 17
$new.Range("B19").Value = "'006143"
$new.Range("C19").Value = "'恒生前海中证质量成长低波动指数A"
$new.Range("D19").Value = "'0.06"
$new.Range("E19").Value = "'94.34"
$new.Range("F19").Value = "'3.87"
$new.Range("G19").Value = "'0.0023"
$new.Range("H19").Value = 1
$new.Range("A20").Value = 18
$new.Range("B20").Value = "'006144"
$new.Range("C20").Value = "'恒生前海中证质量成长低波动指数C"
$new.Range("D20").Value = "'0.01"
$new.Range("E20").Value = "'94.34"
$new.Range("F20").Value = "'3.87"
$new.Range("G20").Value = "'0.0004"
$new.Range("H20").Value = 1
$new.Range("A21").Value = 19
$new.Range("B21").Value = "'014641"
$new.Range("C21").Value = "'上投摩根行业轮动混合C"
$new.Range("D21").Value = "'0.00"
$new.Range("E21").Value = "'86.32"
$new.Range("F21").Value = "'2.44"
$new.Range("G21").Value = 0
$new.Range("H21").Value = 9

# Re-apply the original (unstyled / style index 0) formatting to the data
# block so the cells end up identical to their pre-edit appearance instead
# of picking up the quote-prefix flag from the apostrophe-prefixed input.
$new.Range("B2:G21").Style = $dataStyle

# ---------------------------------------------------------------------------
# Part 2: update "总计" with the new 2022-Q1 summary row
# ---------------------------------------------------------------------------
# Re-resolve "总计" by name: inserting "2022-Q1" shifted its sheet index, and
# the handle captured before the insert still points at the old index (now
# occupied by "2022-Q1"), not at the "总计" worksheet itself.
$zj = $wb.Worksheets.Item("总计")

$zj.Rows.Item(2).Insert()

# The inserted row picks up bold/bordered formatting from the header row
# above it; restore the plain per-column formatting used by the other data
# rows (index column bordered+bold, the rest unstyled) via a format-only
# paste from the row directly below (now the old 2021-Q4 row, row 3).
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)
$zj.Range("B3:D3").Copy()
$zj.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 20
$zj.Range("D2").Value = 4.09

# Renumber the (0-based) index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3
